$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Low-grade glioma")
$ws.Range("B2").Value = 0.158965320748276
$ws.Range("B3").Value = 0.146962372110544
$ws.Range("B4").Value = 0.874193558521807
$ws.Range("B5").Value = 0.728339031091421
$ws.Range("B6").Value = 0.148792161026365
$ws.Range("B7").Value = 0.836515635385441
$ws.Range("B8").Value = 0.196307130417606
$ws.Range("B9").Value = 0.134720453987371
$ws = $wb.Worksheets.Item("Mesenchymal tumor")
$ws.Range("B2").Value = 0.658978328173375
$ws.Range("B3").Value = 0.335913312693499
$ws.Range("B5").Value = 0.201754385964909
$ws.Range("B6").Value = 0.661919504643954
$ws.Range("B7").Value = 0.490842490842491
$ws.Range("B8").Value = 0.851814851814852
$ws.Range("B9").Value = 0.812910413930149
$ws = $wb.Worksheets.Item("Neurofibroma plexiform")
$ws.Range("B9").Value = 0.714399735860193
$ws = $wb.Worksheets.Item("Non-neoplastic tumor")
$ws.Range("B9").Value = 0.371730893209898
$ws = $wb.Worksheets.Item("Germ cell tumor")
$ws.Range("B9").Value = 0.323012594134958
$ws = $wb.Worksheets.Item("Schwannoma")
$ws.Range("B9").Value = 0.0597088054740609
$ws = $wb.Worksheets.Item("Choroid plexus tumor")
$ws.Range("B9").Value = 0.121982838921517
$ws = $wb.Worksheets.Item("Other tumor")
$ws.Range("B2").Value = 0.167077949858817
$ws.Range("B4").Value = 0.383195002994785
$ws.Range("B5").Value = 0.613074356122186
$ws.Range("B6").Value = 0.744072711369706
$ws.Range("B7").Value = 0.750771130026087
$ws.Range("B8").Value = 0.873790272564862
$ws.Range("B9").Value = 0.342330077926341
$ws = $wb.Worksheets.Item("Medulloblastoma")
$ws.Range("B2").Value = 0.178765281671919
$ws.Range("B3").Value = 0.818370221050665
$ws.Range("B4").Value = 0.946079289753066
$ws.Range("B5").Value = 0.971915975844798
$ws.Range("B6").Value = 0.399888854575204
$ws.Range("B7").Value = 0.00307470546413358
$ws.Range("B8").Value = 0.00261898965050277
$ws.Range("B9").Value = 0.0235346934986866
$ws = $wb.Worksheets.Item("Mixed neuronal-glial tumor")
$ws.Range("B2").Value = 0.686514723206952
$ws.Range("B3").Value = 0.0187897276984035
$ws.Range("B4").Value = 0.999999999999969
$ws.Range("B5").Value = 0.498680101486856
$ws.Range("B6").Value = 0.999999999999936
$ws.Range("B7").Value = 0.164184087571725
$ws.Range("B8").Value = 0.206137365005419
$ws.Range("B9").Value = 0.183607803413947
$ws = $wb.Worksheets.Item("Ependymoma")
$ws.Range("B2").Value = 0.106717016062028
$ws.Range("B3").Value = 0.765661101678121
$ws.Range("B4").Value = 0.474483025747333
$ws.Range("B5").Value = 0.293332234357925
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 0.194785106832606
$ws.Range("B8").Value = 0.0979396497170013
$ws.Range("B9").Value = 0.429688254318136
$ws = $wb.Worksheets.Item("Other high-grade glioma")
$ws.Range("B2").Value = 0.19955781429289
$ws.Range("B3").Value = 0.309865431987481
$ws.Range("B4").Value = 0.248632198653321
$ws.Range("B5").Value = 0.0735420797512746
$ws.Range("B6").Value = 0.276011201785736
$ws.Range("B7").Value = 0.155730714345813
$ws.Range("B8").Value = 0.482383043352823
$ws.Range("B9").Value = 0.331349500657068
$ws = $wb.Worksheets.Item("Craniopharyngioma")
$ws.Range("B2").Value = 0.737839053628527
$ws.Range("B3").Value = 0.63088594667542
$ws.Range("B4").Value = 0.751066856330015
$ws.Range("B5").Value = 0.300616405879564
$ws.Range("B6").Value = 0.363411994990942
$ws.Range("B7").Value = 0.668548387096774
$ws.Range("B8").Value = 0.0701612903225806
$ws.Range("B9").Value = 0.519033579392551
$ws = $wb.Worksheets.Item("ATRT")
$ws.Range("B9").Value = 0.150045324913797
$ws = $wb.Worksheets.Item("Meningioma")
$ws.Range("B2").Value = 0.406344451378776
$ws.Range("B3").Value = 0.080976353928299
$ws.Range("B4").Value = 0.480877424126851
$ws.Range("B5").Value = 0.270560062871274
$ws.Range("B6").Value = 0.388985507246377
$ws.Range("B7").Value = 0.914425100183614
$ws.Range("B8").Value = 0.640665835712275
$ws.Range("B9").Value = 0.158298398203477
$ws = $wb.Worksheets.Item("DIPG or DMG")
$ws.Range("B2").Value = 0.285789423282333
$ws.Range("B3").Value = 0.760344776517266
$ws.Range("B4").Value = 0.999999999999994
$ws.Range("B5").Value = 0.870022337929905
$ws.Range("B6").Value = 0.999999999999997
$ws.Range("B7").Value = 0.4606509551112
$ws.Range("B8").Value = 0.291184371184371
$ws.Range("B9").Value = 0.823354276121646
